$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a Text number format on the cells whose new values look numeric
# (e.g. "1.000", "0.7114") so Excel stores them as literal text, matching
# the inline-string cells in the original workbook, instead of parsing them
# into floating point numbers.
$ws.Range("D5,D6,D7,D8,D9,D10,D11,D13,D14,D15,D17,D18,D19,D22,D23,D24,D25,D26,D28,D29,D30,D31,D33,D34,D36,D37,D38,D39,D41,D42,D43,D44,D47,D48,D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.323.79"

$ws.Range("D3").Value = "1.876.81"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "0.7114"
$ws.Range("E5").Value = "  -0.36%  "

$ws.Range("D6").Value = "242.29"
$ws.Range("E6").Value = "  +0.69%  "

$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "0.3111"
$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("D9").Value = "0.07764"
$ws.Range("E9").Value = "  +0.26%  "

$ws.Range("D10").Value = "25.04"
$ws.Range("E10").Value = "  +0.04%  "

$ws.Range("D11").Value = "0.08464"
$ws.Range("E11").Value = "  +2.49%  "

$ws.Range("D12").Value = "1.888.88"
$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("D13").Value = "5.216"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "0.7113"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "91.41"
$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("D16").Value = "29.312.15"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "0.000008290"
$ws.Range("E17").Value = "  +6.14%  "

$ws.Range("D18").Value = "6.008"
$ws.Range("E18").Value = "  +2.32%  "

$ws.Range("D19").Value = "242.77"
$ws.Range("E19").Value = "  -0.71%  "

$ws.Range("D20").Value = "2.133.95"
$ws.Range("E20").Value = "  +1.38%  "

$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "7.844"
$ws.Range("E23").Value = "  -1.31%  "

$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").Value = "0.1612"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").Value = "162.60"
$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D28").Value = "18.50"
$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("D29").Value = "1.514"
$ws.Range("E29").Value = "  +1.24%  "

$ws.Range("D30").Value = "4.404"
$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").Value = "4.334"
$ws.Range("E31").Value = "  +4.28%  "

$ws.Range("E32").Value = "  -2.83%  "

$ws.Range("D33").Value = "0.05240"
$ws.Range("E33").Value = "  +0.82%  "

$ws.Range("D34").Value = "1.933"
$ws.Range("E34").Value = "  +1.25%  "

$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").Value = "0.7413"
$ws.Range("E36").Value = "  +1.73%  "

$ws.Range("D37").Value = "2.688"
$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("D38").Value = "0.01870"
$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").Value = "2.726"
$ws.Range("E39").Value = "  +1.39%  "

$ws.Range("D40").Value = "1.173.50"
$ws.Range("E40").Value = "  +1.52%  "

$ws.Range("D41").Value = "6.390"
$ws.Range("E41").Value = "  +4.76%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8878"
$ws.Range("E42").Value = "  -2.01%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "72.94"
$ws.Range("E43").Value = "  +0.67%  "

$ws.Range("D44").Value = "106.42"
$ws.Range("E44").Value = "  +4.63%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").Value = "2.029.05"
$ws.Range("E46").Value = "  +1.16%  "

$ws.Range("D47").Value = "1.817"
$ws.Range("E47").Value = "  +2.82%  "

$ws.Range("D48").Value = "0.5203"
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("E49").Value = "  +4.41%  "

$ws.Range("D50").Value = "9.406"
$ws.Range("E50").Value = "  +0.78%  "

$ws.Range("E51").Value = "  +1.07%  "

Write-Output "done"